{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"5840560 - Marco Antonio Carvalho Pereira\";\nfor (const para of paragraphs.items) {\n  if (para.text && para.text.trim() === target) {\n    para.delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the bullet paragraph listing the responsible instructor and remove\n# the whole paragraph (text + its paragraph mark), so the heading\n# \"Docente(s) Respons\u00e1vel(eis)\" is immediately followed by the next heading.\n$target = \"5840560 - Marco Antonio Carvalho Pereira\"\n\n$rng = $d.Content\nwhile ($rng.Find.Execute($target)) {\n    $rng.Expand(4) | Out-Null   # wdParagraph - grow to the full paragraph incl. its mark\n    $rng.Delete()\n    $rng = $d.Content\n}\n"}
